# Add a new "UK" worksheet to the workbook, based on the existing "Poland"
# sheet (same layout/styles), fill in the UK-specific values, make it the
# active sheet with cell C17 selected - matching the author's edit that
# added Test Data for the UK Market.

$wb = $excel.ActiveWorkbook

# Poland is the current last sheet; copy it and place the copy right after it.
$poland = $wb.Worksheets.Item("Poland")
$poland.Copy($null, $poland)

# The copy becomes the last sheet in the workbook - rename it to "UK".
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "UK"

# Fill in the UK-specific data (set B4 first so the new shared strings are
# appended in the same order as the target: NGC code, then "UK Market").
$newSheet.Range("B4").Value = "NGC-2741/T3352/T3357"
$newSheet.Range("B2").Value = "UK Market"

# Make the new UK sheet the active tab with C17 selected, mirroring the
# recorded sheet view state.
$newSheet.Activate()
[void]$newSheet.Range("C17").Select()
